# Update cached leve-profit figures (cols H-N) across several sheets to
# reflect refreshed market-board prices, per the scheduled-runner data sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 31228.543
$ws.Range("I64").Value = 202240
$ws.Range("J64").Value = 2726.6333
$ws.Range("K64").Value = 202240
$ws.Range("L64").Value = 2726.6333
$ws.Range("M64").Value = -201992
$ws.Range("N64").Value = -3222.6333

$ws.Range("H67").Value = 31228.543
$ws.Range("I67").Value = 202240
$ws.Range("J67").Value = 2726.6333
$ws.Range("K67").Value = 202240
$ws.Range("L67").Value = 2726.6333
$ws.Range("M67").Value = -201382
$ws.Range("N67").Value = -4442.6333

$ws.Range("H109").Value = 32981.332
$ws.Range("J109").Value = 32981.332
$ws.Range("L109").Value = 32981.332
$ws.Range("N109").Value = -35755.332

$ws.Range("H117").Value = 45743
$ws.Range("J117").Value = 45743
$ws.Range("L117").Value = 45743
$ws.Range("N117").Value = -54921

$ws.Range("H138").Value = 1338.375
$ws.Range("I138").Value = 1024.303
$ws.Range("J138").Value = 2819
$ws.Range("K138").Value = 3072.909000000001
$ws.Range("L138").Value = 8457
$ws.Range("M138").Value = 2067.090999999999
$ws.Range("N138").Value = -18737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2336.6956
$ws.Range("I61").Value = 1113.25
$ws.Range("J61").Value = 3671.3635
$ws.Range("K61").Value = 1113.25
$ws.Range("L61").Value = 3671.3635
$ws.Range("M61").Value = -901.25
$ws.Range("N61").Value = -4095.3635

$ws.Range("H74").Value = 962.5952
$ws.Range("I74").Value = 685.2
$ws.Range("J74").Value = 2349.5715
$ws.Range("K74").Value = 685.2
$ws.Range("L74").Value = 2349.5715
$ws.Range("M74").Value = 188.8
$ws.Range("N74").Value = -4097.5715

$ws.Range("H77").Value = 962.5952
$ws.Range("I77").Value = 685.2
$ws.Range("J77").Value = 2349.5715
$ws.Range("K77").Value = 3426
$ws.Range("L77").Value = 11747.8575
$ws.Range("M77").Value = 942
$ws.Range("N77").Value = -20483.8575

$ws.Range("H80").Value = 49100.285
$ws.Range("J80").Value = 49100.285
$ws.Range("L80").Value = 49100.285
$ws.Range("N80").Value = -51096.285

$ws.Range("H83").Value = 49100.285
$ws.Range("J83").Value = 49100.285
$ws.Range("L83").Value = 147300.855
$ws.Range("N83").Value = -157284.855

$ws.Range("H117").Value = 38557.6
$ws.Range("J117").Value = 38557.6
$ws.Range("L117").Value = 38557.6
$ws.Range("N117").Value = -47735.6

$ws.Range("H136").Value = 2336.6956
$ws.Range("I136").Value = 1113.25
$ws.Range("J136").Value = 3671.3635
$ws.Range("K136").Value = 3339.75
$ws.Range("L136").Value = 11014.0905
$ws.Range("M136").Value = -789.75
$ws.Range("N136").Value = -16114.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 935
$ws.Range("I64").Value = 925
$ws.Range("J64").Value = 945
$ws.Range("K64").Value = 925
$ws.Range("L64").Value = 945
$ws.Range("M64").Value = -700
$ws.Range("N64").Value = -1395

$ws.Range("H67").Value = 935
$ws.Range("I67").Value = 925
$ws.Range("J67").Value = 945
$ws.Range("K67").Value = 925
$ws.Range("L67").Value = 945
$ws.Range("M67").Value = -145
$ws.Range("N67").Value = -2505

$ws.Range("H105").Value = 2185.1794
$ws.Range("I105").Value = 2442.8572
$ws.Range("J105").Value = 2128.8125
$ws.Range("K105").Value = 2442.8572
$ws.Range("L105").Value = 2128.8125
$ws.Range("M105").Value = -695.8571999999999
$ws.Range("N105").Value = -5622.8125

$ws.Range("H117").Value = 46515.668
$ws.Range("J117").Value = 46515.668
$ws.Range("L117").Value = 46515.668
$ws.Range("N117").Value = -55693.668

$ws.Range("H124").Value = 49977.668
$ws.Range("J124").Value = 49977.668
$ws.Range("L124").Value = 49977.668
$ws.Range("N124").Value = -59797.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 24245.72
$ws.Range("I132").Value = 975.64
$ws.Range("J132").Value = 130018.82
$ws.Range("K132").Value = 2926.92
$ws.Range("L132").Value = 390056.46
$ws.Range("M132").Value = -396.9200000000001
$ws.Range("N132").Value = -395116.46

$ws.Range("H134").Value = 182868.33
$ws.Range("I134").Value = 749.7288
$ws.Range("J134").Value = 779812.6
$ws.Range("K134").Value = 2249.1864
$ws.Range("L134").Value = 2339437.8
$ws.Range("M134").Value = 285.8136
$ws.Range("N134").Value = -2344507.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3991
$ws.Range("I113").Value = 9020.083000000001
$ws.Range("J113").Value = 638.2778
$ws.Range("K113").Value = 27060.249
$ws.Range("L113").Value = 1914.8334
$ws.Range("M113").Value = -24890.249
$ws.Range("N113").Value = -6254.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 818
$ws.Range("I22").Value = 747.5
$ws.Range("J22").Value = 874.4
$ws.Range("K22").Value = 747.5
$ws.Range("L22").Value = 874.4
$ws.Range("M22").Value = -452.5
$ws.Range("N22").Value = -1464.4

$ws.Range("H27").Value = 818
$ws.Range("I27").Value = 747.5
$ws.Range("J27").Value = 874.4
$ws.Range("K27").Value = 747.5
$ws.Range("L27").Value = 874.4
$ws.Range("M27").Value = -640.5
$ws.Range("N27").Value = -1088.4

$ws.Range("H81").Value = 20890.8
$ws.Range("J81").Value = 20890.8
$ws.Range("L81").Value = 20890.8
$ws.Range("N81").Value = -22886.8

$ws.Range("H84").Value = 20890.8
$ws.Range("J84").Value = 20890.8
$ws.Range("L84").Value = 62672.39999999999
$ws.Range("N84").Value = -72656.39999999999

$ws.Range("H132").Value = 2054.861
$ws.Range("I132").Value = 1290.3036
$ws.Range("J132").Value = 4730.8125
$ws.Range("K132").Value = 3870.9108
$ws.Range("L132").Value = 14192.4375
$ws.Range("M132").Value = -1340.9108
$ws.Range("N132").Value = -19252.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45881.332
$ws.Range("J16").Value = 45881.332
$ws.Range("L16").Value = 45881.332
$ws.Range("N16").Value = -46465.332

$ws.Range("H132").Value = 1244.5555
$ws.Range("I132").Value = 450.34616
$ws.Range("J132").Value = 3309.5
$ws.Range("K132").Value = 1351.03848
$ws.Range("L132").Value = 9928.5
$ws.Range("M132").Value = 1178.96152
$ws.Range("N132").Value = -14988.5

$ws.Range("H136").Value = 15789.91
$ws.Range("I136").Value = 22769.29
$ws.Range("J136").Value = 1513.909
$ws.Range("K136").Value = 68307.87
$ws.Range("L136").Value = 4541.727000000001
$ws.Range("M136").Value = -65757.87
$ws.Range("N136").Value = -9641.727000000001
